# Adds an "Address" column (new column F) derived from the school/address
# portion of column B, inserting it before the existing "District" column
# (which shifts from F to G). Rows 67 and 68 keep no derived address, matching
# the source data (their original District values could not be cleanly split).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting a whole column before the existing "District" column (F) pushes
# District -> G automatically and leaves a blank column F in its place.
$ws.Columns.Item(6).Insert()

$ws.Range("F2").Value = "Address"

$addresses = @{
    3 = "Govt. High SchoolA S GudiHospet"
    4 = "Govt. P U College(HG) Siruguppa"
    5 = "G H P S MotasugurSiruguppa"
    6 = "Govt. P U College (High School Section) T B Dam Hospet"
    7 = "G H P S HaregondanahalliH B Halli"
    8 = "G H P S SangameshwaraH B Halli"
    9 = "Shree Bikki Naresh BabuGovt. High SchoolKudithni"
    10 = "B R C Siruguppa"
    11 = "G H P S H HosahalliSiruguppa"
    12 = "Govt. Girls High SchoolKurugodu"
    13 = "S G P U (Comp) CollegeY Nagghashastri Nagar"
    14 = "G H P SAnekal ThandaH B Halli"
    15 = "G H P SGududur"
    16 = "G H S Sidiginamola"
    17 = "Sri Sidalighewra High SchoolSovenahalli Sandur"
    18 = "V D H L High SchoolS G P Road"
    19 = "G H S Nehru ColonyHosapete"
    20 = "G H S Kallukambha"
    21 = "G H SBadanahatti"
    22 = "Govt. High School"
    23 = "G H S PujarahallyKudligi"
    24 = "G H S Karur SchoolSiruguppa"
    25 = "G P G High SchoolHuvina Hadagalli"
    26 = "G H S IttigiHadagali"
    27 = "S G Comp P U CollegeY Nagesh Shastri Nagar"
    28 = "Netaji Subhash Chandra Bose Govt. Urdu High SchoolMillerpet"
    29 = "Govt Juniar CollegeBommanahalli"
    30 = "K V H High School KattebenurHadagali"
    31 = "G H S TimalapuraKudligi"
    32 = "Govt. P U CollegeRadiopark"
    33 = "G H P SchoolVaradapurH B Halli"
    34 = "P M C Girls High SchoolHospet"
    35 = "Sri Renuka High SchoolH B Halli"
    36 = "G P U C High School SectionChoranurSandur"
    37 = "Vijaya High School Kampli KSF Hospet"
    38 = "Huttina Yellamma G H S RaraviSiruguppa"
    39 = "G H P SThimmalapur"
    40 = "Govt. Urdu High SchoolHiriyalkudamCowl BazarEd Gah Road"
    41 = "G H P S Dasarahalli TandaHuvina Hadagali"
    42 = "Govt. High SchoolB M SugurSiruguppa"
    43 = "Govt. Adarsha VidyalayaSandur"
    44 = "Govt. High School Yelubenchi"
    45 = "G H P S Yelubenchi"
    46 = "G H P S MaduruH B Halli"
    47 = "G G H STekkalakoteSiraguppa"
    48 = "G H P S AgraharaSandur"
    49 = "B R R G High SchoolG NagalapurHospet"
    50 = "L H S High SchoolHampasagarHagaribommanahalli"
    51 = "G J CollegeEmmiganuru"
    52 = "Govt. High SchoolByasigideriHagaribommanahalli"
    53 = "G P U C Hosamoka"
    54 = "Govt. High SchoolM SugurSiraguppa"
    55 = "S G B H SBennikalliH B Halli"
    56 = "G H S ToranagalluSandur"
    57 = "G H S Siddammanahalli"
    58 = "G H S Donimali Sandur"
    59 = "G G H S Hosahalli Kudligi"
    60 = "G P U CollegeChittawadgiHospet"
    61 = "G H P S GosubaluSirugappa"
    62 = "Smiode vyasapuriHigh School Vyasa Nakare Hospet"
    63 = "S B G H S MeerakoranahalliHuvinahadagali"
    64 = "Govt. High School RavihalSiraguppa"
    65 = "G H P S B N Halli"
    66 = "Govt. High SchoolSirigeriSiruguppa"
    69 = "G H P S K Gudda Siruguppa"
}

foreach ($row in $addresses.Keys) {
    $ws.Range("F$row").Value = $addresses[$row]
}

